$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper idiom used throughout: toggling a character-level formatting
# property on and back off for a tiny Range forces the engine to keep a run
# boundary at that point instead of silently re-merging text with identical
# neighbouring runs. It leaves no visible formatting residue behind.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 1) Insert a new list paragraph right after "... access the signing Url."
#    paragraph, before the first blank paragraph that precedes "Pros:".
#    The paragraph naturally inherits the ListParagraph / numId=1 list
#    formatting from paragraph 19 (itself a numbered list item), matching
#    the target pPr (pStyle=ListParagraph, numPr ilvl=0 numId=1).
# ---------------------------------------------------------------------------
$signingUrlPara = $d.Paragraphs(19)
$signingUrlPara.Range.InsertParagraphAfter()

$newListPara = $d.Paragraphs(20)
$newListPara.Range.Text = "After user signed the agreement, we can call the View Agreement API endpoint to view the signed document (just open the Url under DOUCMNENT type name, not MANAGE). "

# Split "Url" into its own run (mirrors the proofErr-wrapped run in the
# target, minus the spell-check marker itself which the object model does
# not expose).
$prefix = "After user signed the agreement, we can call the View Agreement API endpoint to view the signed document (just open the "
$urlStart = $newListPara.Range.Start + $prefix.Length
$urlEnd = $urlStart + 3
$urlRange = $d.Range($urlStart, $urlEnd)
$urlRange.Font.Bold = 1
$urlRange.Font.Bold = 0
$tailPoint = $d.Range($urlEnd, $urlEnd)
$tailPoint.Font.Italic = 1
$tailPoint.Font.Italic = 0

Write-Output "Step 1 done"

# ---------------------------------------------------------------------------
# 2) Collapse the three runs of the "Adobe eSignature offer cheaper pricing
#    plan ..." sentence into a single run (no text change, only run
#    consolidation). Find/Execute over the whole sentence forces the engine
#    to re-emit it as one run.
# ---------------------------------------------------------------------------
$pricingText = "Adobe eSignature offer cheaper pricing plan which required USD 14.99 monthly per business license."
$d.Content.Find.Execute($pricingText, $true, $false, $false, $false, $false, $true, 1, $false, $pricingText, 2) | Out-Null

Write-Output "Step 2 done"

# ---------------------------------------------------------------------------
# 3) Isolate "email." into its own run (the lastRenderedPageBreak that used
#    to sit before "For authentication..." is relocated in front of it in
#    real Word; the object model has no settable property for that pure
#    rendering/pagination marker, so only the run split is reproduced here).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("user email.", $true, $false, $false, $false, $false, $true, 1, $false, "user email.", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("email.", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r.Font.Bold = 1
$r.Font.Bold = 0

Write-Output "Step 3 done"

# ---------------------------------------------------------------------------
# 4) "OAuth 2.0 will be problematic" -> "OAuth will be difficult", with
#    "difficult" split into its own run (matching the target's 3-run split).
#    This edit also removes the lastRenderedPageBreak that used to sit at
#    the start of this run, matching the target (it no longer appears here).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("OAuth 2.0 will be problematic", $true, $false, $false, $false, $false, $true, 1, $false, "OAuth will be difficult", 2) | Out-Null

$r2 = $d.Content
$r2.Find.Execute("difficult", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r2.Font.Bold = 1
$r2.Font.Bold = 0

Write-Output "Step 4 done"

# ---------------------------------------------------------------------------
# 5) Split the hyperlink display text into 3 runs (no text change):
#    "https://community.adobe.com/t5/adobe" + "-" +
#    "acrobat-sign-discussions/help-with-oauth2/td-p/7063211"
# ---------------------------------------------------------------------------
$hlRange = $d.Content
$hlRange.Find.Execute("https://community.adobe.com/t5/adobe-acrobat-sign-discussions/help-with-oauth2/td-p/7063211", $true) | Out-Null
$hlStart = $hlRange.Start
$dashStart = $hlStart + 36
$dashEnd = $dashStart + 1
$dashRange = $d.Range($dashStart, $dashEnd)
$dashRange.Font.Bold = 1
$dashRange.Font.Bold = 0
$afterDash = $d.Range($dashEnd, $dashEnd)
$afterDash.Font.Bold = 1
$afterDash.Font.Bold = 0

Write-Output "Step 5 done"

# ---------------------------------------------------------------------------
# 6) styles.xml: FollowedHyperlink character style needs to be materialised
#    (it previously only existed as a latent-style exception entry). Apply
#    it to a throw-away one-character range so the engine emits the full
#    <w:style> element, then tune the priority/visibility flags to match
#    real Word's default template for this built-in style as closely as
#    the object model allows.
# ---------------------------------------------------------------------------
$followedHyperlinkStyle = $d.Styles(-87)
$styleProbe = $d.Range($hlStart, $hlStart + 1)
$savedText = $styleProbe.Text
$styleProbe.Style = $followedHyperlinkStyle
$styleProbe.Style = "Hyperlink"
$fhStyle = $d.Styles(-87)
$fhStyle.Priority = 99
$fhStyle.UnhideWhenUsed = $true

Write-Output "Step 6 done"
